$d = $word.ActiveDocument

# Step 1: extend the existing "Draws 1 card..." run with the new trailing
# text via Find/Replace so it inherits the exact same run formatting
# (sz/szCs/lang) as the text it is appended to.
$old = "Draws 1 card for each draw modifier on the card"
$new = "Draws 1 card for each draw modifier on the card when played."
$rng = $d.Content
$rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# Step 2: locate the newly appended tail (" when played.") and give it its
# own (no-underline) run formatting. Word will automatically split the run
# at the formatting boundary, producing two new runs: " " and "when played."
$tailRng = $d.Content
$tailRng.Find.Execute(" when played.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$spaceRng = $d.Range($tailRng.Start, $tailRng.Start + 1)
$spaceRng.Font.Underline = 0

$restRng = $d.Range($tailRng.Start + 1, $tailRng.End)
$restRng.Font.Underline = 0
